$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "is_locked" (E) and "is_enabled" (F) columns from the usr import
# template header row. Everything to the right (dept_ids_lbl, role_ids_lbl,
# rem) shifts left, and the now-unused trailing cells disappear.
$ws.Range("E1:F1").EntireColumn.Delete()
